$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Source Code" styled paragraphs that followed
#    "About this site" (the rendered R console input/output block).
# ------------------------------------------------------------------
$startPos = -1
$endPos = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Source Code") {
        if ($startPos -eq -1) {
            $startPos = $p.Range.Start
        }
        $endPos = $p.Range.End
    }
}
if ($startPos -ne -1) {
    $killRange = $d.Range($startPos, $endPos)
    $killRange.Delete()
}

# ------------------------------------------------------------------
# 2. Add the new "Abstract Title" paragraph style, placed (in intent)
#    right before the existing "Abstract" style and used as the style
#    preceding it.
# ------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# ------------------------------------------------------------------
# 3. Tighten the spacing above the "Abstract" style (before: 300 -> 100)
# ------------------------------------------------------------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ------------------------------------------------------------------
# 4. Add the new "Footnote Block Text" paragraph style, based on and
#    following "Footnote Text".
# ------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
